$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 439.40475
$ws.Cells.Item(135, 9).Value = 305.2903
$ws.Cells.Item(135, 10).Value = 817.36365
$ws.Cells.Item(135, 11).Value = 2747.6127
$ws.Cells.Item(135, 12).Value = 7356.27285
$ws.Cells.Item(135, 13).Value = -212.6127000000001
$ws.Cells.Item(135, 14).Value = -12426.27285

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10531296
$ws.Cells.Item(32, 9).Value = 12198056
$ws.Cells.Item(32, 10).Value = 17884.46
$ws.Cells.Item(32, 11).Value = 12198056
$ws.Cells.Item(32, 12).Value = 17884.46
$ws.Cells.Item(32, 13).Value = -12197769
$ws.Cells.Item(32, 14).Value = -18458.46
$ws.Cells.Item(132, 8).Value = 1602.6731
$ws.Cells.Item(132, 9).Value = 1342.2858
$ws.Cells.Item(132, 10).Value = 2696.3
$ws.Cells.Item(132, 11).Value = 4026.8574
$ws.Cells.Item(132, 12).Value = 8088.900000000001
$ws.Cells.Item(132, 13).Value = -1496.8574
$ws.Cells.Item(132, 14).Value = -13148.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1773.973
$ws.Cells.Item(134, 9).Value = 1555.9678
$ws.Cells.Item(134, 10).Value = 2900.3333
$ws.Cells.Item(134, 11).Value = 4667.903399999999
$ws.Cells.Item(134, 12).Value = 8700.999899999999
$ws.Cells.Item(134, 13).Value = -2132.903399999999
$ws.Cells.Item(134, 14).Value = -13770.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2050.5469
$ws.Cells.Item(31, 9).Value = 1150.3125
$ws.Cells.Item(31, 10).Value = 4751.25
$ws.Cells.Item(31, 11).Value = 1150.3125
$ws.Cells.Item(31, 12).Value = 4751.25
$ws.Cells.Item(31, 13).Value = -855.3125
$ws.Cells.Item(31, 14).Value = -5341.25
$ws.Cells.Item(34, 8).Value = 2050.5469
$ws.Cells.Item(34, 9).Value = 1150.3125
$ws.Cells.Item(34, 10).Value = 4751.25
$ws.Cells.Item(34, 11).Value = 1150.3125
$ws.Cells.Item(34, 12).Value = 4751.25
$ws.Cells.Item(34, 13).Value = -948.3125
$ws.Cells.Item(34, 14).Value = -5155.25
$ws.Cells.Item(58, 8).Value = 949.1348
$ws.Cells.Item(58, 9).Value = 556.8305
$ws.Cells.Item(58, 11).Value = 556.8305
$ws.Cells.Item(58, 13).Value = -353.8305
$ws.Cells.Item(74, 8).Value = 15010
$ws.Cells.Item(74, 10).Value = 15834.875
$ws.Cells.Item(74, 12).Value = 15834.875
$ws.Cells.Item(74, 14).Value = -17582.875
$ws.Cells.Item(77, 8).Value = 15010
$ws.Cells.Item(77, 10).Value = 15834.875
$ws.Cells.Item(77, 12).Value = 47504.625
$ws.Cells.Item(77, 14).Value = -56240.625
$ws.Cells.Item(107, 8).Value = 601.6957
$ws.Cells.Item(107, 9).Value = 405.45456
$ws.Cells.Item(107, 10).Value = 781.5833
$ws.Cells.Item(107, 11).Value = 405.45456
$ws.Cells.Item(107, 12).Value = 781.5833
$ws.Cells.Item(107, 13).Value = 1514.54544
$ws.Cells.Item(107, 14).Value = -4621.5833
$ws.Cells.Item(132, 8).Value = 1852.0193
$ws.Cells.Item(132, 9).Value = 1946.3429
$ws.Cells.Item(132, 10).Value = 1657.8235
$ws.Cells.Item(132, 11).Value = 5839.028700000001
$ws.Cells.Item(132, 12).Value = 4973.470499999999
$ws.Cells.Item(132, 13).Value = -3309.028700000001
$ws.Cells.Item(132, 14).Value = -10033.4705
$ws.Cells.Item(134, 8).Value = 1740.8136
$ws.Cells.Item(134, 9).Value = 1123.9788
$ws.Cells.Item(134, 10).Value = 4156.75
$ws.Cells.Item(134, 11).Value = 3371.936400000001
$ws.Cells.Item(134, 12).Value = 12470.25
$ws.Cells.Item(134, 13).Value = -836.9364000000005
$ws.Cells.Item(134, 14).Value = -17540.25
$ws.Cells.Item(136, 8).Value = 949.1348
$ws.Cells.Item(136, 9).Value = 556.8305
$ws.Cells.Item(136, 11).Value = 1670.4915
$ws.Cells.Item(136, 13).Value = 879.5084999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 837.73334
$ws.Cells.Item(5, 9).Value = 733.75
$ws.Cells.Item(5, 10).Value = 956.5714
$ws.Cells.Item(5, 11).Value = 2201.25
$ws.Cells.Item(5, 12).Value = 2869.7142
$ws.Cells.Item(5, 13).Value = -2089.25
$ws.Cells.Item(5, 14).Value = -3093.7142
$ws.Cells.Item(60, 8).Value = 321
$ws.Cells.Item(60, 9).Value = 150
$ws.Cells.Item(60, 10).Value = 342.375
$ws.Cells.Item(60, 11).Value = 450
$ws.Cells.Item(60, 12).Value = 1027.125
$ws.Cells.Item(60, 13).Value = -199
$ws.Cells.Item(60, 14).Value = -1529.125
$ws.Cells.Item(75, 8).Value = 1565.1428
$ws.Cells.Item(75, 9).Value = 713
$ws.Cells.Item(75, 10).Value = 1906
$ws.Cells.Item(75, 11).Value = 2139
$ws.Cells.Item(75, 12).Value = 5718
$ws.Cells.Item(75, 13).Value = -1141
$ws.Cells.Item(75, 14).Value = -7714
$ws.Cells.Item(78, 8).Value = 1565.1428
$ws.Cells.Item(78, 9).Value = 713
$ws.Cells.Item(78, 10).Value = 1906
$ws.Cells.Item(78, 11).Value = 6417
$ws.Cells.Item(78, 12).Value = 17154
$ws.Cells.Item(78, 13).Value = -1425
$ws.Cells.Item(78, 14).Value = -27138
$ws.Cells.Item(132, 8).Value = 29462028
$ws.Cells.Item(132, 9).Value = 1161.3334
$ws.Cells.Item(132, 10).Value = 58922896
$ws.Cells.Item(132, 11).Value = 10452.0006
$ws.Cells.Item(132, 12).Value = 530306064
$ws.Cells.Item(132, 13).Value = -7922.000599999999
$ws.Cells.Item(132, 14).Value = -530311124
$ws.Cells.Item(134, 8).Value = 3340.3845
$ws.Cells.Item(134, 9).Value = 1060.7142
$ws.Cells.Item(134, 11).Value = 3182.1426
$ws.Cells.Item(134, 13).Value = 1887.8574
$ws.Cells.Item(135, 8).Value = 837.73334
$ws.Cells.Item(135, 9).Value = 733.75
$ws.Cells.Item(135, 10).Value = 956.5714
$ws.Cells.Item(135, 11).Value = 6603.75
$ws.Cells.Item(135, 12).Value = 8609.142600000001
$ws.Cells.Item(135, 13).Value = -4068.75
$ws.Cells.Item(135, 14).Value = -13679.1426
$ws.Cells.Item(136, 8).Value = 3572.4443
$ws.Cells.Item(136, 9).Value = 2917
$ws.Cells.Item(136, 10).Value = 5866.5
$ws.Cells.Item(136, 11).Value = 8751
$ws.Cells.Item(136, 12).Value = 17599.5
$ws.Cells.Item(136, 13).Value = -3651
$ws.Cells.Item(136, 14).Value = -27799.5
$ws.Cells.Item(137, 8).Value = 4374.077
$ws.Cells.Item(137, 9).Value = 3985.7144
$ws.Cells.Item(137, 10).Value = 4827.1665
$ws.Cells.Item(137, 11).Value = 11957.1432
$ws.Cells.Item(137, 12).Value = 14481.4995
$ws.Cells.Item(137, 13).Value = -6857.143199999999
$ws.Cells.Item(137, 14).Value = -24681.4995
$ws.Cells.Item(138, 8).Value = 1642
$ws.Cells.Item(138, 9).Value = 797.9
$ws.Cells.Item(138, 10).Value = 2291.3076
$ws.Cells.Item(138, 11).Value = 2393.7
$ws.Cells.Item(138, 12).Value = 6873.9228
$ws.Cells.Item(138, 13).Value = 2746.3
$ws.Cells.Item(138, 14).Value = -17153.9228
$ws.Cells.Item(139, 8).Value = 3614.2104
$ws.Cells.Item(139, 9).Value = 1147.7778
$ws.Cells.Item(139, 10).Value = 5834
$ws.Cells.Item(139, 11).Value = 3443.3334
$ws.Cells.Item(139, 12).Value = 17502
$ws.Cells.Item(139, 13).Value = 1696.6666
$ws.Cells.Item(139, 14).Value = -27782
$ws.Cells.Item(140, 8).Value = 4766257.5
$ws.Cells.Item(140, 9).Value = 6252982
$ws.Cells.Item(140, 10).Value = 8739.799999999999
$ws.Cells.Item(140, 11).Value = 18758946
$ws.Cells.Item(140, 12).Value = 26219.4
$ws.Cells.Item(140, 13).Value = -18753766
$ws.Cells.Item(140, 14).Value = -36579.39999999999
$ws.Cells.Item(141, 8).Value = 3701.8235
$ws.Cells.Item(141, 9).Value = 1559
$ws.Cells.Item(141, 10).Value = 6112.5
$ws.Cells.Item(141, 11).Value = 4677
$ws.Cells.Item(141, 12).Value = 18337.5
$ws.Cells.Item(141, 13).Value = 503
$ws.Cells.Item(141, 14).Value = -28697.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 6650
$ws.Cells.Item(9, 9).Value = 900
$ws.Cells.Item(9, 10).Value = 12400
$ws.Cells.Item(9, 11).Value = 900
$ws.Cells.Item(9, 12).Value = 12400
$ws.Cells.Item(9, 13).Value = -730
$ws.Cells.Item(9, 14).Value = -12740

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(87, 8).Value = 49989
$ws.Cells.Item(87, 10).Value = 49989
$ws.Cells.Item(87, 12).Value = 49989
$ws.Cells.Item(87, 14).Value = -52235
$ws.Cells.Item(90, 8).Value = 49989
$ws.Cells.Item(90, 10).Value = 49989
$ws.Cells.Item(90, 12).Value = 149967
$ws.Cells.Item(90, 14).Value = -161199
$ws.Cells.Item(132, 8).Value = 1379.97
$ws.Cells.Item(132, 9).Value = 1412.4674
$ws.Cells.Item(132, 10).Value = 1006.25
$ws.Cells.Item(132, 11).Value = 4237.4022
$ws.Cells.Item(132, 12).Value = 3018.75
$ws.Cells.Item(132, 13).Value = -1707.4022
$ws.Cells.Item(132, 14).Value = -8078.75
$ws.Cells.Item(133, 8).Value = 45142.95
$ws.Cells.Item(133, 10).Value = 45142.95
$ws.Cells.Item(133, 12).Value = 45142.95
$ws.Cells.Item(133, 14).Value = -50202.95
$ws.Cells.Item(136, 8).Value = 2506
$ws.Cells.Item(136, 9).Value = 2012.175
$ws.Cells.Item(136, 10).Value = 4152.0835
$ws.Cells.Item(136, 11).Value = 6036.525
$ws.Cells.Item(136, 12).Value = 12456.2505
$ws.Cells.Item(136, 13).Value = -3486.525
$ws.Cells.Item(136, 14).Value = -17556.2505

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 7053.778
$ws.Cells.Item(96, 9).Value = 1308.4
$ws.Cells.Item(96, 10).Value = 9263.538
$ws.Cells.Item(96, 11).Value = 1308.4
$ws.Cells.Item(96, 12).Value = 9263.538
$ws.Cells.Item(96, 13).Value = 64.59999999999991
$ws.Cells.Item(96, 14).Value = -12009.538
